# Electricity sector updates to:
# 1) fix subscript ordering in MPCbS (Capacity Cost Multiplier vs Share of
#    Existing Capacity Built this Year curve, raised to the 1.3 power);
# 2) update offshore wind capacity / resource calibration params
#    (Share of Cost Effective Capacity Built in a Single Year);
# 3) modify electricity parameters to better calibrate the power sector.

$wb = $excel.ActiveWorkbook

# Remember the sheet that should remain active/selected when we're done
# (the workbook keeps "About" as its active tab throughout this edit).
$wsAbout = $wb.Worksheets.Item("About")

# --- Sheet: CSC-CSCCCMvSoECBtY --------------------------------------------
# "CSC Capacity Supply Curve Capacity Cost Multiplier vs Share of Existing
# Capacity Built this Year" - update the cost-multiplier curve (row 2,
# columns D:O) to the recalibrated values (old_value ^ 1.3).
$wsMult = $wb.Worksheets.Item("CSC-CSCCCMvSoECBtY")

$wsMult.Range("D2").Value = 1.5045088484257014
$wsMult.Range("E2").Value = 2.0640090320749205
$wsMult.Range("F2").Value = 2.8315774206270898
$wsMult.Range("G2").Value = 3.8845908928057571
$wsMult.Range("H2").Value = 5.3292014149580273
$wsMult.Range("I2").Value = 7.3110369981623169
$wsMult.Range("J2").Value = 10.029882116058715
$wsMult.Range("K2").Value = 13.759817560663983
$wsMult.Range("L2").Value = 18.87684990968436
$wsMult.Range("M2").Value = 25.896815910390654
$wsMult.Range("N2").Value = 35.527382916586511
$wsMult.Range("O2").Value = 48.739387133679415

# Update the sheet's remembered selection to match the edited workbook.
$wsMult.Range("D3:Q4").Select()

# --- Sheet: CSC-CSCSoCECBiaSY ---------------------------------------------
# "CSC Capacity Supply Curve Share of Cost Effective Capacity Built in a
# Single Year" - recalibrate the share-built-per-year parameters: most
# resources (hard coal, natural gas steam turbine, natural gas combined
# cycle, nuclear, hydro, solar thermal, biomass, geothermal, municipal
# solid waste, hard coal w CCS, natural gas combined cycle w CCS, biomass
# w CCS, lignite w CCS, small modular reactor, hydrogen combustion
# turbine, hydrogen combined cycle) move from 0.33 to 0.25, while onshore
# wind moves from 0.04 to 0.08. (Rows for solar PV, crude oil, and heavy
# or residual fuel oil were already 0.25/0/0 and are unchanged.)
$wsShare = $wb.Worksheets.Item("CSC-CSCSoCECBiaSY")

$wsShare.Range("B2:AE6").Value = 0.25
$wsShare.Range("B7:AE7").Value = 0.08
$wsShare.Range("B9:AE15").Value = 0.25
$wsShare.Range("B18:AE25").Value = 0.25

# Update the sheet's remembered selection to match the edited workbook.
$wsShare.Range("I22").Select()

# Restore "About" as the active sheet/tab (selecting ranges above on other
# sheets switches the active sheet as a side effect).
$wsAbout.Activate()
